# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.390.96'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '3.262.40'
$ws.Range('E3').Value = '  +2.59%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''614.19'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').Value = '''157.74'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.261.22'
$ws.Range('E8').Value = '  +2.54%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +1.59%  '
$ws.Range('D11').Value = '''5.79'
$ws.Range('E11').Value = '  +2.00%  '
$ws.Range('E12').Value = '  -4.31%  '
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = '''39.09'
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('D15').Value = '3.798.72'
$ws.Range('E15').Value = '  +2.72%  '
$ws.Range('D16').Value = '66.433.77'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').Value = '''7.44'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '3.264.32'
$ws.Range('E18').Value = '  +2.81%  '
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('D20').Value = '''505.10'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').Value = '''15.46'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = '''0.755'
$ws.Range('E22').Value = '  +3.14%  '
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').Value = '''14.65'
$ws.Range('E24').Value = '  -1.46%  '
$ws.Range('D25').Value = '''87.21'
$ws.Range('E25').Value = '  +2.97%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('D28').Value = '''9.23'
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '''0.134'
$ws.Range('E29').Value = '  +51.83%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '''2.39'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').Value = '''7.12'
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('E32').Value = '  -4.09%  '
$ws.Range('D33').Value = '''27.98'
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('E35').Value = '  -4.05%  '
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('E37').Value = '  +18.90%  '
$ws.Range('D38').Value = '''55.80'
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('D39').Value = '0.0₃0788'
$ws.Range('E39').Value = '  +15.09%  '
$ws.Range('D40').Value = '''497.35'
$ws.Range('E40').Value = '  -1.74%  '
$ws.Range('D41').Value = '''0.0423'
$ws.Range('E41').Value = '  +0.52%  '
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('D43').Value = '''8.84'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '''2.53'
$ws.Range('E44').Value = '  +3.54%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '''0.294'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '3.007.53'
$ws.Range('E46').Value = '  +6.42%  '
$ws.Range('D47').Value = '''29.01'
$ws.Range('E47').Value = '  +3.09%  '
$ws.Range('E48').Value = '  +5.43%  '
$ws.Range('E49').Value = '  +2.16%  '
$ws.Range('E51').Value = '  -3.52%  '
